$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D16").Value = 4
$ws.Range("D17").Value = 3
$ws.Range("D18").Value = 3

$lo = $ws.ListObjects.Item(2)
$lo.Resize($ws.Range("B2:G18"))
$ws.Range("H2").Value = ""

$ws.Range("M19").Select()
